$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 268, pushing existing rows 268-335 down to 270-337
$ws.Rows.Item(268).Resize(2).Insert()

# Row 268: new data (Primera)
$ws.Cells.Item(268, 1).Value = 11
$ws.Cells.Item(268, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(268, 3).Value = "Bíobío"
$ws.Cells.Item(268, 4).Value = 44964
$ws.Cells.Item(268, 5).Value = 8
$ws.Cells.Item(268, 6).Value = 100114013
$ws.Cells.Item(268, 7).Value = "Zanahoria"
$ws.Cells.Item(268, 8).Value = "Sin especificar"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 800
$ws.Cells.Item(268, 11).Value = 8000
$ws.Cells.Item(268, 12).Value = 8500
$ws.Cells.Item(268, 13).Value = 8250
$ws.Cells.Item(268, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(268, 15).Value = "Región de Ñuble"
$ws.Cells.Item(268, 16).Value = 412
$ws.Cells.Item(268, 17).Value = 20
$ws.Cells.Item(268, 18).Value = "Hortaliza"

# Row 269: new data (Segunda)
$ws.Cells.Item(269, 1).Value = 11
$ws.Cells.Item(269, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(269, 3).Value = "Bíobío"
$ws.Cells.Item(269, 4).Value = 44964
$ws.Cells.Item(269, 5).Value = 8
$ws.Cells.Item(269, 6).Value = 100114013
$ws.Cells.Item(269, 7).Value = "Zanahoria"
$ws.Cells.Item(269, 8).Value = "Sin especificar"
$ws.Cells.Item(269, 9).Value = "Segunda"
$ws.Cells.Item(269, 10).Value = 400
$ws.Cells.Item(269, 11).Value = 7000
$ws.Cells.Item(269, 12).Value = 7000
$ws.Cells.Item(269, 13).Value = 7000
$ws.Cells.Item(269, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(269, 15).Value = "Región de Ñuble"
$ws.Cells.Item(269, 16).Value = 350
$ws.Cells.Item(269, 17).Value = 20
$ws.Cells.Item(269, 18).Value = "Hortaliza"

# Apply the same date-cell number format as other D column cells
$ws.Cells.Item(268, 4).NumberFormat = $ws.Cells.Item(267, 4).NumberFormat
$ws.Cells.Item(269, 4).NumberFormat = $ws.Cells.Item(267, 4).NumberFormat
